$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column in the middle of the uniform-width C:D block so the
# inserted column inherits the 15.25 width exactly (avoids lossy rounding
# that happens when setting ColumnWidth directly). This pushes old D -> E
# and old E ("wait") -> F.
$ws.Columns("D").Insert()

# The insert left the original "storeTitle"/json-value column (old D) now
# sitting at E, and a blank (but correctly-widthed) column at D. Move the
# original D content back to D, leaving E as the new duplicate column.
$ws.Range("E1:E3").Copy()
$ws.Range("D1:D3").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# E2 should stay empty/absent (same as before the edit for the old E2 cell
# position), so drop whatever got copied into it.
$ws.Range("E2").ClearContents()

# E3 becomes the new "title2" stored-value cell (mirrors D2's JSON value).
$ws.Range("E3").Value = '{"value":"title2"}'

# Update the active selection to C6.
$ws.Range("C6").Select()
